$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.218.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.340.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.80%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.90%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.81%  '

# Row 7
$ws.Range("E7").Value = '  -0.18%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.30%  '

# Row 9
$ws.Range("E9").Value = '  +5.34%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.583'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.76%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.56%  '

# Row 12
$ws.Range("E12").Value = '  +2.53%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '692.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.95%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.884.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.66%  '

# Row 15
$ws.Range("E15").Value = '  +1.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.304.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.06%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.119'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.36%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.339.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.61%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.89%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.895'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.25%  '

# Row 22
$ws.Range("E22").Value = '  +1.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.84%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.75%  '

# Row 25
$ws.Range("E25").Value = '  +2.36%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.59%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.18%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.57%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.94%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.14%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '566.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.54%  '

# Row 32
$ws.Range("E32").Value = '  +1.56%  '

# Row 33
$ws.Range("E33").Value = '  +2.28%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.57%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.710.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.22%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.53%  '

# Row 38
$ws.Range("E38").Value = '  +3.96%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.96%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.58%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.86%  '

# Row 42
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.336'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.96%  '

# Row 43
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0673'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.64%  '

# Row 44
$ws.Range("E44").Value = '  +0.59%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0414'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.40%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.84%  '

# Row 47
$ws.Range("E47").Value = '  +1.45%  '

# Row 48
$ws.Range("E48").Value = '  -0.29%  '

# Row 49
$ws.Range("E49").Value = '  -0.07%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.98%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.74%  '
